$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.915.10'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '3.271.48'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '182.54'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.598'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.67'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('D12').Value = '3.836.96'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.58'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '68.838.52'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000172'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').Value = '3.270.93'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.84'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.58'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '394.30'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.96%  '
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.97'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.517'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  +5.55%  '
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '23.01'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.15'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '164.02'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.95'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.830'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '26.41'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.60'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.63'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '41.42'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.49'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.96%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0690'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '345.74'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.20%  '
$ws.Range('D46').Value = '2.610.44'
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '24.79'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.32'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '31.71'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('E51').Value = '  -0.31%  '
